$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row for new columns
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"
$ws.Range("G1").Value = "Terms Typically Offered"

# Row 2
$ws.Range("D2").Value = "NA"
$ws.Range("E2").Value = "NA"
$ws.Range("F2").Value = "NA"
$ws.Range("G2").Value = "F"

# Row 3
$ws.Range("C3").Value = "GRC 101; Graphic Communication Majors and Minors only."
$ws.Range("D3").Value = "NA"
$ws.Range("E3").Value = "NA"
$ws.Range("F3").Value = "NA"
$ws.Range("G3").Value = "W, SP"

# Row 4
$ws.Range("D4").Value = "NA"
$ws.Range("E4").Value = "NA"
$ws.Range("F4").Value = "NA"
$ws.Range("G4").Value = "TBD"

# Row 5
$ws.Range("D5").Value = "NA"
$ws.Range("E5").Value = "NA"
$ws.Range("F5").Value = "NA"
$ws.Range("G5").Value = "F, W"

# Row 6
$ws.Range("D6").Value = "NA"
$ws.Range("E6").Value = "NA"
$ws.Range("F6").Value = "NA"
$ws.Range("G6").Value = "W, SP"

# Row 7
$ws.Range("D7").Value = "NA"
$ws.Range("E7").Value = "NA"
$ws.Range("F7").Value = "NA"
$ws.Range("G7").Value = "W, SU"

# Row 8
$ws.Range("D8").Value = "NA"
$ws.Range("E8").Value = "NA"
$ws.Range("F8").Value = "NA"
$ws.Range("G8").Value = "F, W"

# Row 9
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("F9").Value = "NA"
$ws.Range("G9").Value = "F, W"

# Row 10
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("F10").Value = "NA"
$ws.Range("G10").Value = "F, W, SP"

# Row 11
$ws.Range("D11").Value = "NA"
$ws.Range("E11").Value = "NA"
$ws.Range("F11").Value = "NA"
$ws.Range("G11").Value = "W, SP"

# Row 12
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "NA"
$ws.Range("G12").Value = "TBD"

# Row 13
$ws.Range("D13").Value = "NA"
$ws.Range("E13").Value = "NA"
$ws.Range("F13").Value = "NA"
$ws.Range("G13").Value = "F, SP"

# Row 14
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "NA"
$ws.Range("F14").Value = "NA"
$ws.Range("G14").Value = "F, W"

# Row 15
$ws.Range("D15").Value = "NA"
$ws.Range("E15").Value = "NA"
$ws.Range("F15").Value = "NA"
$ws.Range("G15").Value = "F, W"

# Row 16
$ws.Range("D16").Value = "NA"
$ws.Range("E16").Value = "NA"
$ws.Range("F16").Value = "NA"
$ws.Range("G16").Value = "F, W"

# Row 17
$ws.Range("D17").Value = "NA"
$ws.Range("E17").Value = "NA"
$ws.Range("F17").Value = "NA"
$ws.Range("G17").Value = "F"

# Row 18
$ws.Range("D18").Value = "NA"
$ws.Range("E18").Value = "NA"
$ws.Range("F18").Value = "NA"
$ws.Range("G18").Value = "F, SP"

# Row 19
$ws.Range("D19").Value = "NA"
$ws.Range("E19").Value = "NA"
$ws.Range("F19").Value = "NA"
$ws.Range("G19").Value = "F, SP"

# Row 20
$ws.Range("D20").Value = "NA"
$ws.Range("E20").Value = "NA"
$ws.Range("F20").Value = "NA"
$ws.Range("G20").Value = "W"

# Row 21
$ws.Range("D21").Value = "NA"
$ws.Range("E21").Value = "NA"
$ws.Range("F21").Value = "NA"
$ws.Range("G21").Value = "F, W"

# Row 22
$ws.Range("D22").Value = "NA"
$ws.Range("E22").Value = "NA"
$ws.Range("F22").Value = "NA"
$ws.Range("G22").Value = "W, SP"

# Row 23
$ws.Range("D23").Value = "NA"
$ws.Range("E23").Value = "NA"
$ws.Range("F23").Value = "NA"
$ws.Range("G23").Value = "W, SP"

# Row 24
$ws.Range("D24").Value = "NA"
$ws.Range("E24").Value = "NA"
$ws.Range("F24").Value = "NA"
$ws.Range("G24").Value = "F"

# Row 25
$ws.Range("D25").Value = "NA"
$ws.Range("E25").Value = "NA"
$ws.Range("F25").Value = "NA"
$ws.Range("G25").Value = "W"

# Row 26
$ws.Range("D26").Value = "NA"
$ws.Range("E26").Value = "NA"
$ws.Range("F26").Value = "NA"
$ws.Range("G26").Value = "F, SP"

# Row 27
$ws.Range("D27").Value = "NA"
$ws.Range("E27").Value = "NA"
$ws.Range("F27").Value = "NA"
$ws.Range("G27").Value = "F, W, SP"

# Row 28
$ws.Range("D28").Value = "NA"
$ws.Range("E28").Value = "NA"
$ws.Range("F28").Value = "NA"
$ws.Range("G28").Value = "F, W, SP"

# Row 29
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("F29").Value = "NA"
$ws.Range("G29").Value = "TBD"

# Row 30
$ws.Range("D30").Value = "NA"
$ws.Range("E30").Value = "NA"
$ws.Range("F30").Value = "NA"
$ws.Range("G30").Value = "F, SP"

# Row 31
$ws.Range("D31").Value = "NA"
$ws.Range("E31").Value = "NA"
$ws.Range("F31").Value = "NA"
$ws.Range("G31").Value = "F, SP"

# Row 32
$ws.Range("D32").Value = "NA"
$ws.Range("E32").Value = "NA"
$ws.Range("F32").Value = "NA"
$ws.Range("G32").Value = "W, SP"

# Row 33
$ws.Range("C33").Value = "Senior standing; JOUR 331; JOUR 342; BUS 453; and one of the GRC 338, GRC 377, or JOUR 390."
$ws.Range("D33").Value = "NA"
$ws.Range("E33").Value = "NA"
$ws.Range("F33").Value = "NA"
$ws.Range("G33").Value = "SP"

# Row 34
$ws.Range("D34").Value = "NA"
$ws.Range("E34").Value = "NA"
$ws.Range("F34").Value = "NA"
$ws.Range("G34").Value = "SP"

# Row 35
$ws.Range("D35").Value = "NA"
$ws.Range("E35").Value = "NA"
$ws.Range("F35").Value = "NA"
$ws.Range("G35").Value = "F, W"

# Row 36
$ws.Range("D36").Value = "NA"
$ws.Range("E36").Value = "NA"
$ws.Range("F36").Value = "NA"
$ws.Range("G36").Value = "W"

# Row 37
$ws.Range("D37").Value = "NA"
$ws.Range("E37").Value = "NA"
$ws.Range("F37").Value = "NA"
$ws.Range("G37").Value = "W"

# Row 38
$ws.Range("C38").Value = "GRC 337; Graphic Communication major."
$ws.Range("D38").Value = "NA"
$ws.Range("E38").Value = "NA"
$ws.Range("F38").Value = "GRC 361."
$ws.Range("G38").Value = "TBD "

# Row 39
$ws.Range("D39").Value = "NA"
$ws.Range("E39").Value = "NA"
$ws.Range("F39").Value = "NA"
$ws.Range("G39").Value = "W"

# Row 40
$ws.Range("D40").Value = "NA"
$ws.Range("E40").Value = "NA"
$ws.Range("F40").Value = "NA"
$ws.Range("G40").Value = "SP"

# Row 41
$ws.Range("D41").Value = "NA"
$ws.Range("E41").Value = "NA"
$ws.Range("F41").Value = "NA"
$ws.Range("G41").Value = "TBD"

# Row 42
$ws.Range("D42").Value = "NA"
$ws.Range("E42").Value = "NA"
$ws.Range("F42").Value = "NA"
$ws.Range("G42").Value = "TBD"

# Row 43
$ws.Range("D43").Value = "NA"
$ws.Range("E43").Value = "NA"
$ws.Range("F43").Value = "NA"
$ws.Range("G43").Value = "TBD"

# Row 44
$ws.Range("D44").Value = "NA"
$ws.Range("E44").Value = "NA"
$ws.Range("F44").Value = "NA"
$ws.Range("G44").Value = "TBD"

# Row 45
$ws.Range("D45").Value = "NA"
$ws.Range("E45").Value = "NA"
$ws.Range("F45").Value = "NA"
$ws.Range("G45").Value = "F, W, SP"

# Row 46
$ws.Range("D46").Value = "NA"
$ws.Range("E46").Value = "NA"
$ws.Range("F46").Value = "NA"
$ws.Range("G46").Value = "TBD"

# Row 47
$ws.Range("D47").Value = "NA"
$ws.Range("E47").Value = "NA"
$ws.Range("F47").Value = "NA"
$ws.Range("G47").Value = "TBD"

# Row 48
$ws.Range("D48").Value = "NA"
$ws.Range("E48").Value = "NA"
$ws.Range("F48").Value = "NA"
$ws.Range("G48").Value = "TBD"

# Row 49
$ws.Range("D49").Value = "NA"
$ws.Range("E49").Value = "NA"
$ws.Range("F49").Value = "NA"
$ws.Range("G49").Value = "F,W,SP,SU"

# Row 50
$ws.Range("D50").Value = "NA"
$ws.Range("E50").Value = "NA"
$ws.Range("F50").Value = "NA"
$ws.Range("G50").Value = "F,W,SP,SU"

# Row 51
$ws.Range("D51").Value = "NA"
$ws.Range("E51").Value = "NA"
$ws.Range("F51").Value = "NA"
$ws.Range("G51").Value = "TBD"

# Row 52
$ws.Range("D52").Value = "NA"
$ws.Range("E52").Value = "NA"
$ws.Range("F52").Value = "NA"
$ws.Range("G52").Value = "TBD"

# Row 53
$ws.Range("D53").Value = "NA"
$ws.Range("E53").Value = "NA"
$ws.Range("F53").Value = "NA"
$ws.Range("G53").Value = "TBD"

# Row 54
$ws.Range("D54").Value = "GRC 501."
$ws.Range("E54").Value = "NA"
$ws.Range("F54").Value = "NA"
$ws.Range("G54").Value = "TBD"

# Row 55
$ws.Range("D55").Value = "GRC 501."
$ws.Range("E55").Value = "NA"
$ws.Range("F55").Value = "NA"
$ws.Range("G55").Value = "TBD"

# Row 56
$ws.Range("D56").Value = "GRC 501."
$ws.Range("E56").Value = "NA"
$ws.Range("F56").Value = "NA"
$ws.Range("G56").Value = "TBD"

# Row 57
$ws.Range("D57").Value = "NA"
$ws.Range("E57").Value = "NA"
$ws.Range("F57").Value = "NA"
$ws.Range("G57").Value = "TBD"

# Row 58
$ws.Range("D58").Value = "NA"
$ws.Range("E58").Value = "NA"
$ws.Range("F58").Value = "NA"
$ws.Range("G58").Value = "TBD"

# Row 59
$ws.Range("D59").Value = "GRC 502."
$ws.Range("E59").Value = "NA"
$ws.Range("F59").Value = "NA"
$ws.Range("G59").Value = "TBD"

# Row 60
$ws.Range("D60").Value = "GRC 502."
$ws.Range("E60").Value = "NA"
$ws.Range("F60").Value = "NA"
$ws.Range("G60").Value = "TBD"

# Row 61
$ws.Range("D61").Value = "GRC 502."
$ws.Range("E61").Value = "NA"
$ws.Range("F61").Value = "NA"
$ws.Range("G61").Value = "TBD"

# Row 62
$ws.Range("D62").Value = "GRC 530."
$ws.Range("E62").Value = "NA"
$ws.Range("F62").Value = "NA"
$ws.Range("G62").Value = "TBD"

# Row 63
$ws.Range("D63").Value = "NA"
$ws.Range("E63").Value = "NA"
$ws.Range("F63").Value = "NA"
$ws.Range("G63").Value = "TBD"

# Row 64
$ws.Range("D64").Value = "NA"
$ws.Range("E64").Value = "NA"
$ws.Range("F64").Value = "NA"
$ws.Range("G64").Value = "TBD"

